$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "22.077.73"
Set-TextValue $ws.Range("E2") "  -1.47%  "
Set-TextValue $ws.Range("D3") "1.557.65"
Set-TextValue $ws.Range("E3") "  -0.56%  "
Set-TextValue $ws.Range("D4") "0.9992"
Set-TextValue $ws.Range("E4") "  -0.18%  "
Set-TextValue $ws.Range("D5") "0.9992"
Set-TextValue $ws.Range("E5") "  -0.18%  "
Set-TextValue $ws.Range("D6") "287.40"
Set-TextValue $ws.Range("E6") "  +0.40%  "
Set-TextValue $ws.Range("D7") "0.3842"
Set-TextValue $ws.Range("E7") "  +3.67%  "
Set-TextValue $ws.Range("D8") "0.3245"
Set-TextValue $ws.Range("E8") "  -1.19%  "
Set-TextValue $ws.Range("D9") "41.52"
Set-TextValue $ws.Range("E9") "  -10.77%  "
Set-TextValue $ws.Range("D10") "1.127"
Set-TextValue $ws.Range("E10") "  -1.79%  "
Set-TextValue $ws.Range("E11") "  -1.08%  "
Set-TextValue $ws.Range("D12") "0.9993"
Set-TextValue $ws.Range("E12") "  -0.19%  "
Set-TextValue $ws.Range("D13") "19.44"
Set-TextValue $ws.Range("E13") "  -5.12%  "
Set-TextValue $ws.Range("D14") "5.717"
Set-TextValue $ws.Range("E14") "  -2.19%  "
Set-TextValue $ws.Range("D15") "6.814"
Set-TextValue $ws.Range("E15") "  +0.03%  "
Set-TextValue $ws.Range("D16") "1.556.47"
Set-TextValue $ws.Range("E16") "  -0.52%  "
Set-TextValue $ws.Range("D17") "0.00001098"
Set-TextValue $ws.Range("E17") "  -0.31%  "
Set-TextValue $ws.Range("D18") "0.06621"
Set-TextValue $ws.Range("E18") "  -1.28%  "
Set-TextValue $ws.Range("D19") "85.17"
Set-TextValue $ws.Range("E19") "  -1.25%  "
Set-TextValue $ws.Range("D20") "6.410"
Set-TextValue $ws.Range("E20") "  +1.23%  "
Set-TextValue $ws.Range("D22") "15.99"
Set-TextValue $ws.Range("E22") "  -1.89%  "
Set-TextValue $ws.Range("E23") "  -2.64%  "
Set-TextValue $ws.Range("D24") "22.088.13"
Set-TextValue $ws.Range("E24") "  -1.38%  "
Set-TextValue $ws.Range("D25") "2.332"
Set-TextValue $ws.Range("E25") "  +0.47%  "
Set-TextValue $ws.Range("D26") "2.555"
Set-TextValue $ws.Range("E26") "  -0.69%  "
Set-TextValue $ws.Range("D27") "148.97"
Set-TextValue $ws.Range("E27") "  -1.19%  "
Set-TextValue $ws.Range("D28") "18.88"
Set-TextValue $ws.Range("E28") "  -2.63%  "
Set-TextValue $ws.Range("D29") "4.859"
Set-TextValue $ws.Range("E29") "  -1.89%  "
Set-TextValue $ws.Range("D30") "1.730.99"
Set-TextValue $ws.Range("E30") "  -0.44%  "
Set-TextValue $ws.Range("D31") "120.74"
Set-TextValue $ws.Range("E31") "  -2.39%  "
Set-TextValue $ws.Range("D32") "1.113"
Set-TextValue $ws.Range("E32") "  +5.63%  "
Set-TextValue $ws.Range("D33") "5.892"
Set-TextValue $ws.Range("E33") "  -1.58%  "
Set-TextValue $ws.Range("D34") "1.686"
Set-TextValue $ws.Range("E34") "  -14.59%  "
Set-TextValue $ws.Range("D35") "9.308"
Set-TextValue $ws.Range("E35") "  -3.96%  "
Set-TextValue $ws.Range("D36") "0.08183"
Set-TextValue $ws.Range("E36") "  -0.97%  "
Set-TextValue $ws.Range("D37") "0.06229"
Set-TextValue $ws.Range("E37") "  -1.52%  "
Set-TextValue $ws.Range("D38") "0.02298"
Set-TextValue $ws.Range("E38") "  -4.34%  "
Set-TextValue $ws.Range("D39") "5.220"
Set-TextValue $ws.Range("E39") "  -0.02%  "
Set-TextValue $ws.Range("D40") "0.2109"
Set-TextValue $ws.Range("E40") "  -3.56%  "
Set-TextValue $ws.Range("D41") "1.223"
Set-TextValue $ws.Range("E41") "  -6.33%  "
Set-TextValue $ws.Range("D42") "10.92"
Set-TextValue $ws.Range("E42") "  -2.10%  "
Set-TextValue $ws.Range("D43") "0.9984"
Set-TextValue $ws.Range("E43") "  -0.23%  "
Set-TextValue $ws.Range("D44") "0.5967"
Set-TextValue $ws.Range("E44") "  -2.55%  "
Set-TextValue $ws.Range("D45") "13.54"
Set-TextValue $ws.Range("E45") "  -0.95%  "
Set-TextValue $ws.Range("D46") "3.718"
Set-TextValue $ws.Range("E46") "  -0.79%  "
Set-TextValue $ws.Range("D47") "0.5762"
Set-TextValue $ws.Range("E47") "  -3.26%  "
Set-TextValue $ws.Range("D48") "1.936"
Set-TextValue $ws.Range("E48") "  -3.92%  "
Set-TextValue $ws.Range("D49") "119.43"
Set-TextValue $ws.Range("E49") "  -3.65%  "
Set-TextValue $ws.Range("E50") "  -2.02%  "
Set-TextValue $ws.Range("D51") "0.06898"
Set-TextValue $ws.Range("E51") "  -3.66%  "
